$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row right after the current data (row 38 in this workbook)
$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

# Force the new cells to be stored as text (matching the existing sheet's
# convention of storing every value, numeric-looking or not, as a shared
# string) so Excel does not auto-convert the date / price into a real
# date serial or number.
$newRange = $ws.Range("A" + $newRow + ":D" + $newRow)
$newRange.NumberFormat = "@"

$ws.Range("A" + $newRow).Value = "2026-02-07"
$ws.Range("B" + $newRow).Value = "6510380"
$ws.Range("C" + $newRow).Value = "0"
$ws.Range("D" + $newRow).Value = "0"

# Drop the explicit number-format style again so the new row ends up using
# the same default (unstyled) cell formatting as every other row.
$newRange.Style = "Normal"
